$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.621.94"
$ws.Range("E2").Value = "  -3.84%  "
$ws.Range("D3").Value = "2.917.60"
$ws.Range("E3").Value = "  -2.19%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "547.65"
$ws.Range("E5").Value = "  -3.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.76"
$ws.Range("E6").Value = "  +4.60%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.509"
$ws.Range("E8").Value = "  +2.04%  "
$ws.Range("D9").Value = "2.910.11"
$ws.Range("E9").Value = "  -2.19%  "
$ws.Range("E10").Value = "  -3.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.72"
$ws.Range("E11").Value = "  -4.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.443"
$ws.Range("E12").Value = "  +0.75%  "
$ws.Range("E13").Value = "  -0.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.64"
$ws.Range("E14").Value = "  +1.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.123"
$ws.Range("E15").Value = "  +1.84%  "
$ws.Range("D16").Value = "3.399.26"
$ws.Range("E16").Value = "  -2.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.88"
$ws.Range("E17").Value = "  +6.51%  "
$ws.Range("D18").Value = "2.913.06"
$ws.Range("E18").Value = "  -2.53%  "
$ws.Range("D19").Value = "57.603.74"
$ws.Range("E19").Value = "  -4.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "416.40"
$ws.Range("E20").Value = "  -2.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.24"
$ws.Range("E21").Value = "  +1.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.689"
$ws.Range("E22").Value = "  +3.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.22"
$ws.Range("E23").Value = "  +3.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.96"
$ws.Range("E24").Value = "  -0.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.44"
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("E28").Value = "  -2.08%  "
$ws.Range("E29").Value = "  +3.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.36"
$ws.Range("E30").Value = "  +3.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.17"
$ws.Range("E31").Value = "  +0.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.94"
$ws.Range("E32").Value = "  -2.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0963"
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.67"
$ws.Range("E34").Value = "  +2.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.931"
$ws.Range("E35").Value = "  +1.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.05"
$ws.Range("E36").Value = "  +4.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "48.11"
$ws.Range("E37").Value = "  -4.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.74"
$ws.Range("E38").Value = "  +3.19%  "
$ws.Range("D39").Value = "0.0₃0684"
$ws.Range("E39").Value = "  +5.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.55"
$ws.Range("E40").Value = "  +5.91%  "
$ws.Range("E41").Value = "  +0.65%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "2.703.52"
$ws.Range("E42").Value = "  +1.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0343"
$ws.Range("E43").Value = "  -2.59%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "371.66"
$ws.Range("E44").Value = "  +1.03%  "
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "123.69"
$ws.Range("E46").Value = "  +2.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.234"
$ws.Range("E47").Value = "  +1.19%  "
$ws.Range("E48").Value = "  +0.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.94"
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.78"
$ws.Range("E50").Value = "  -1.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.98"
$ws.Range("E51").Value = "  -0.32%  "
